# Update performance dashboard 2025-12-20 10:06
#
# Refreshes the Pattern1-Pure Data metrics (deepseek-v3, gemini-3-pro, gpt-5,
# llama-3.1-405b) on both the "Summary" sheet (rows 2-5) and the
# "Pattern1-Pure Data" sheet (rows 2-5), which mirror each other.
#
# Several of the target strings look like numbers/percentages/dates
# ("+0.45%", "20251219", ...). Assigning such text straight to .Value makes
# Excel's COM layer "smart type" them into actual numbers. To keep them as
# literal text (matching the source workbook's inline strings) without
# leaving a stray NumberFormat/quote-prefix style behind on the cell, each
# value is staged in a scratch cell that is explicitly Text-formatted, then
# copied across with Paste Special (values only) and the scratch cell is
# cleared again.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($ws, [string]$addr, [string]$value)

    $scratch = $ws.Range("Z100")
    $scratch.NumberFormat = "@"
    $scratch.Value = $value
    $scratch.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
    $scratch.Clear()
}

$sheetNames = @("Summary", "Pattern1-Pure Data")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Sheets.Item($sheetName)

    # --- Row 2: deepseek-v3 ---
    Set-TextValue $ws "D2" "¥1,004,526.00"
    Set-TextValue $ws "E2" "¥+4,526.00"
    Set-TextValue $ws "F2" "+0.45%"
    Set-TextValue $ws "G2" "+76.65%"
    $ws.Range("H2").Value = 19.872
    Set-TextValue $ws "I2" "0.00%"
    Set-TextValue $ws "J2" "100.0%"
    Set-TextValue $ws "K2" "0.2262%"
    Set-TextValue $ws "L2" "0.1807%"

    # --- Row 3: gemini-3-pro ---
    Set-TextValue $ws "D3" "¥1,004,601.00"
    Set-TextValue $ws "E3" "¥+4,601.00"
    Set-TextValue $ws "F3" "+0.46%"
    Set-TextValue $ws "G3" "+78.32%"
    $ws.Range("H3").Value = 28.141
    Set-TextValue $ws "K3" "0.2299%"
    Set-TextValue $ws "L3" "0.1297%"
    $ws.Range("M3").Value = 3
    Set-TextValue $ws "O3" "20251219"

    # --- Row 4: gpt-5 ---
    Set-TextValue $ws "D4" "¥1,003,469.00"
    Set-TextValue $ws "E4" "¥+3,469.00"
    Set-TextValue $ws "F4" "+0.35%"
    Set-TextValue $ws "G4" "+54.70%"
    $ws.Range("H4").Value = 21.573
    Set-TextValue $ws "I4" "0.00%"
    Set-TextValue $ws "J4" "100.0%"
    Set-TextValue $ws "K4" "0.1734%"
    Set-TextValue $ws "L4" "0.1276%"

    # --- Row 5: llama-3.1-405b ---
    Set-TextValue $ws "D5" "¥1,001,074.00"
    Set-TextValue $ws "E5" "¥+1,074.00"
    Set-TextValue $ws "F5" "+0.11%"
    Set-TextValue $ws "G5" "+14.48%"
    $ws.Range("H5").Value = 9.238
    Set-TextValue $ws "I5" "0.04%"
    Set-TextValue $ws "J5" "50.0%"
    Set-TextValue $ws "K5" "0.0537%"
    Set-TextValue $ws "L5" "0.0923%"
}
